$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value2 = 2.96
$ws.Range("G2").Value2 = 4.2
$ws.Range("I2").Value2 = 2.86
$ws.Range("K2").Value2 = 3.5
$ws.Range("S2").Value2 = 1.01
$ws.Range("T2").Value2 = 2.3
$ws.Range("V2").Value2 = 1.53
$ws.Range("K3").Value2 = 11
$ws.Range("T3").Value2 = 1.98
$ws.Range("U3").Value2 = 1.83
$ws.Range("L4").Value2 = 1.54
$ws.Range("M4").Value2 = 1.12
$ws.Range("P4").Value2 = 1.53
$ws.Range("Q4").Value2 = 2.54
$ws.Range("S4").Value2 = 5.1
$ws.Range("T4").Value2 = 1.98
$ws.Range("U4").Value2 = 1.47
$ws.Range("W4").Value2 = 1.99
$ws.Range("AB4").Value2 = 6.6
$ws.Range("U5").Value2 = 1.96
$ws.Range("X5").Value2 = 9.6
$ws.Range("AB5").Value2 = 8.6
$ws.Range("AM5").Value2 = 140
$ws.Range("F6").Value2 = 1.61
$ws.Range("G6").Value2 = 1.67
$ws.Range("J6").Value2 = 5
$ws.Range("K6").Value2 = 5.6
$ws.Range("P6").Value2 = 3.3
$ws.Range("R6").Value2 = 1.92
$ws.Range("S6").Value2 = 1.88
$ws.Range("T6").Value2 = 1.46
$ws.Range("U6").Value2 = 2.74
$ws.Range("W6").Value2 = 2.48
$ws.Range("AD6").Value2 = 21
$ws.Range("AF6").Value2 = 16.5
$ws.Range("AG6").Value2 = 11
$ws.Range("AI6").Value2 = 44
$ws.Range("AJ6").Value2 = 19.5
$ws.Range("AO6").Value2 = 29
$ws.Range("J7").Value2 = 3.7
$ws.Range("K7").Value2 = 4
$ws.Range("L7").Value2 = 1.4
$ws.Range("U7").Value2 = 1.96
$ws.Range("V7").Value2 = 2.1
$ws.Range("Y7").Value2 = 10
$ws.Range("Z7").Value2 = 13
$ws.Range("AC7").Value2 = 10
$ws.Range("AD7").Value2 = 11
$ws.Range("AL7").Value2 = 80
$ws.Range("J8").Value2 = 3.4
$ws.Range("Q8").Value2 = 2.36
$ws.Range("S8").Value2 = 4.2
$ws.Range("X8").Value2 = 12
$ws.Range("AN9").Value2 = 8.6
$ws.Range("AO9").Value2 = 360
$ws.Range("F10").Value2 = 1.61
$ws.Range("G10").Value2 = 1.62
$ws.Range("H10").Value2 = 6.2
$ws.Range("I10").Value2 = 6.4
$ws.Range("J10").Value2 = 4.5
$ws.Range("K10").Value2 = 4.6
$ws.Range("O10").Value2 = 1.28
$ws.Range("P10").Value2 = 2.16
$ws.Range("Q10").Value2 = 1.84
$ws.Range("U10").Value2 = 2.08
$ws.Range("V10").Value2 = 1.18
$ws.Range("W10").Value2 = 2.6
$ws.Range("AE10").Value2 = 85
$ws.Range("AG10").Value2 = 9.6
$ws.Range("AO10").Value2 = 90
$ws.Range("G11").Value2 = 2.32
$ws.Range("J11").Value2 = 3.1
$ws.Range("K11").Value2 = 3.5
$ws.Range("N11").Value2 = 2.88
$ws.Range("O11").Value2 = 1.44
$ws.Range("P11").Value2 = 1.64
$ws.Range("Q11").Value2 = 2.3
$ws.Range("T11").Value2 = 1.94
$ws.Range("V11").Value2 = 1.3
$ws.Range("W11").Value2 = 1.75
